$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -1219
$ws.Range("H51").Value = 2695.238
$ws.Range("I51").Value = 1390
$ws.Range("J51").Value = 3347.8572
$ws.Range("K51").Value = 1390
$ws.Range("L51").Value = 3347.8572
$ws.Range("M51").Value = -906
$ws.Range("N51").Value = -4315.8572
$ws.Range("H125").Value = 612.3125
$ws.Range("I125").Value = 449.36365
$ws.Range("J125").Value = 970.8
$ws.Range("K125").Value = 4044.27285
$ws.Range("L125").Value = 8737.199999999999
$ws.Range("M125").Value = -1584.27285
$ws.Range("N125").Value = -13657.2
$ws.Range("H137").Value = 2943011.5
$ws.Range("I137").Value = 5001565.5
$ws.Range("K137").Value = 15004696.5
$ws.Range("M137").Value = -15002146.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4996196.5
$ws.Range("I32").Value = 5669593.5
$ws.Range("J32").Value = 29892.625
$ws.Range("K32").Value = 5669593.5
$ws.Range("L32").Value = 29892.625
$ws.Range("M32").Value = -5669306.5
$ws.Range("N32").Value = -30466.625
$ws.Range("H61").Value = 47715300
$ws.Range("I61").Value = 77001144
$ws.Range("J61").Value = 125803.5
$ws.Range("K61").Value = 77001144
$ws.Range("L61").Value = 125803.5
$ws.Range("M61").Value = -77000932
$ws.Range("N61").Value = -126227.5
$ws.Range("H122").Value = 4631680
$ws.Range("I122").Value = 2068.5789
$ws.Range("J122").Value = 22224202
$ws.Range("K122").Value = 6205.736699999999
$ws.Range("L122").Value = 66672606
$ws.Range("M122").Value = -3755.736699999999
$ws.Range("N122").Value = -66677506
$ws.Range("H136").Value = 47715300
$ws.Range("I136").Value = 77001144
$ws.Range("J136").Value = 125803.5
$ws.Range("K136").Value = 231003432
$ws.Range("L136").Value = 377410.5
$ws.Range("M136").Value = -231000882
$ws.Range("N136").Value = -382510.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 633.42
$ws.Range("I31").Value = 396.03845
$ws.Range("J31").Value = 716.82434
$ws.Range("K31").Value = 396.03845
$ws.Range("L31").Value = 716.82434
$ws.Range("M31").Value = -101.03845
$ws.Range("N31").Value = -1306.82434
$ws.Range("H34").Value = 633.42
$ws.Range("I34").Value = 396.03845
$ws.Range("J34").Value = 716.82434
$ws.Range("K34").Value = 396.03845
$ws.Range("L34").Value = 716.82434
$ws.Range("M34").Value = -194.03845
$ws.Range("N34").Value = -1120.82434
$ws.Range("H99").Value = 2293.5217
$ws.Range("I99").Value = 2410.05
$ws.Range("J99").Value = 1516.6666
$ws.Range("K99").Value = 2410.05
$ws.Range("L99").Value = 1516.6666
$ws.Range("M99").Value = -912.0500000000002
$ws.Range("N99").Value = -4512.6666
$ws.Range("H107").Value = 1352.7894
$ws.Range("I107").Value = 799
$ws.Range("J107").Value = 1968.1111
$ws.Range("K107").Value = 799
$ws.Range("L107").Value = 1968.1111
$ws.Range("M107").Value = 1121
$ws.Range("N107").Value = -5808.1111
$ws.Range("H126").Value = 2293.5217
$ws.Range("I126").Value = 2410.05
$ws.Range("J126").Value = 1516.6666
$ws.Range("K126").Value = 7230.150000000001
$ws.Range("L126").Value = 4549.9998
$ws.Range("M126").Value = -4760.150000000001
$ws.Range("N126").Value = -9489.9998
$ws.Range("H132").Value = 46266.914
$ws.Range("I132").Value = 2727.5625
$ws.Range("J132").Value = 145785.42
$ws.Range("K132").Value = 8182.6875
$ws.Range("L132").Value = 437356.26
$ws.Range("M132").Value = -5652.6875
$ws.Range("N132").Value = -442416.26

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1066.069
$ws.Range("I68").Value = 557.3333
$ws.Range("J68").Value = 1376.963
$ws.Range("K68").Value = 1671.9999
$ws.Range("L68").Value = 4130.889
$ws.Range("M68").Value = -860.9999
$ws.Range("N68").Value = -5752.889
$ws.Range("H71").Value = 1066.069
$ws.Range("I71").Value = 557.3333
$ws.Range("J71").Value = 1376.963
$ws.Range("K71").Value = 5015.9997
$ws.Range("L71").Value = 12392.667
$ws.Range("M71").Value = -959.9997000000003
$ws.Range("N71").Value = -20504.667
$ws.Range("H102").Value = 3700
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H107").Value = 1095.6714
$ws.Range("I107").Value = 505.07693
$ws.Range("J107").Value = 1838.6774
$ws.Range("K107").Value = 1515.23079
$ws.Range("L107").Value = 5516.0322
$ws.Range("M107").Value = 404.7692099999999
$ws.Range("N107").Value = -9356.0322
$ws.Range("H131").Value = 761.08
$ws.Range("I131").Value = 509.9
$ws.Range("J131").Value = 928.5333000000001
$ws.Range("K131").Value = 1529.7
$ws.Range("L131").Value = 2785.5999
$ws.Range("M131").Value = 3510.3
$ws.Range("N131").Value = -12865.5999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2553.3333
$ws.Range("I126").Value = 1730
$ws.Range("K126").Value = 5190
$ws.Range("M126").Value = -2720
$ws.Range("H132").Value = 46719.58
$ws.Range("I132").Value = 34131.484
$ws.Range("J132").Value = 74593.21000000001
$ws.Range("K132").Value = 102394.452
$ws.Range("L132").Value = 223779.63
$ws.Range("M132").Value = -99864.45199999999
$ws.Range("N132").Value = -228839.63

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1754.5454
$ws.Range("I68").Value = 1687.5
$ws.Range("J68").Value = 1933.3334
$ws.Range("K68").Value = 1687.5
$ws.Range("L68").Value = 1933.3334
$ws.Range("M68").Value = -938.5
$ws.Range("N68").Value = -3431.3334
$ws.Range("H71").Value = 1754.5454
$ws.Range("I71").Value = 1687.5
$ws.Range("J71").Value = 1933.3334
$ws.Range("K71").Value = 8437.5
$ws.Range("L71").Value = 9666.666999999999
$ws.Range("M71").Value = -4693.5
$ws.Range("N71").Value = -17154.667
$ws.Range("H122").Value = 4421
$ws.Range("I122").Value = 3800
$ws.Range("J122").Value = 4835
$ws.Range("K122").Value = 11400
$ws.Range("L122").Value = 14505
$ws.Range("M122").Value = -8950
$ws.Range("N122").Value = -19405
$ws.Range("H132").Value = 114755.78
$ws.Range("I132").Value = 2001.3334
$ws.Range("J132").Value = 171133
$ws.Range("K132").Value = 6004.0002
$ws.Range("L132").Value = 513399
$ws.Range("M132").Value = -3474.0002
$ws.Range("N132").Value = -518459
$ws.Range("H133").Value = 36029.332
$ws.Range("J133").Value = 36029.332
$ws.Range("L133").Value = 36029.332
$ws.Range("N133").Value = -41089.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 492
$ws.Range("I107").Value = 492
$ws.Range("K107").Value = 1476
$ws.Range("M107").Value = 444
$ws.Range("H126").Value = 1201.4166
$ws.Range("I126").Value = 878.5294
$ws.Range("J126").Value = 1985.5714
$ws.Range("K126").Value = 2635.5882
$ws.Range("L126").Value = 5956.7142
$ws.Range("M126").Value = -165.5882000000001
$ws.Range("N126").Value = -10896.7142
$ws.Range("H132").Value = 50518.78
$ws.Range("I132").Value = 35079.434
$ws.Range("J132").Value = 92626.09
$ws.Range("K132").Value = 105238.302
$ws.Range("L132").Value = 277878.27
$ws.Range("M132").Value = -102708.302
$ws.Range("N132").Value = -282938.27
$ws.Range("H136").Value = 41816.8
$ws.Range("I136").Value = 27311.342
$ws.Range("K136").Value = 81934.026
$ws.Range("M136").Value = -79384.026
